$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62 (shifts the former row 62 -> row 63,
# leaves row 61's original values in place for now).
$ws.Rows.Item(62).Insert()

# The newly inserted row 62 should receive a copy of what row 61 used to
# contain (the "Primera" / Región Metropolitana record dated 44641).
$ws.Cells.Item(62, 1).Value = 6
$ws.Cells.Item(62, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44641
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100102
$ws.Cells.Item(62, 8).Value = "Cítricos"
$ws.Cells.Item(62, 9).Value = 100102006
$ws.Cells.Item(62, 10).Value = "Pomelo"
$ws.Cells.Item(62, 11).Value = "Start Ruby"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 16
$ws.Cells.Item(62, 14).Value = 180000
$ws.Cells.Item(62, 15).Value = 180000
$ws.Cells.Item(62, 16).Value = 180000
$ws.Cells.Item(62, 17).Value = "$/bins (350 kilos)"
$ws.Cells.Item(62, 18).Value = "Región Metropolitana"
$ws.Cells.Item(62, 19).Value = 514
$ws.Cells.Item(62, 20).Value = 350

# Row 61 now becomes the updated record (new date/volume/prices/origin).
$ws.Cells.Item(61, 4).Value = 45223
$ws.Cells.Item(61, 13).Value = 20
$ws.Cells.Item(61, 14).Value = 150000
$ws.Cells.Item(61, 15).Value = 150000
$ws.Cells.Item(61, 16).Value = 150000
$ws.Cells.Item(61, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(61, 19).Value = 429
